# Update the lattice-multiplication practice table: every one of the 15
# cells gets a new "AxB" problem (new operands, new partial-product grid).
# The table shape (5 rows x 3 cols) and each cell's run/line layout stay
# the same - only the digits inside change.
#
# We rebuild each cell's single paragraph via Range.InsertXML with a
# literal OOXML <w:p> fragment (rather than Find/Replace on the text)
# so the "  C    D" / "  ----" lines keep their xml:space="preserve"
# attribute exactly as the target markup expects.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellContent($row, $col, $runXml) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range.Paragraphs.Item(1).Range
    $pkg = '<?xml version="1.0"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
             '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
                 '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body><w:p>' + $runXml + '</w:p></w:body>' +
                 '</w:document>' +
               '</pkg:xmlData>' +
             '</pkg:part>' +
           '</pkg:package>'
    $rng.InsertXML($pkg)
}

function Cell-Run($line1, $line2, $line3, $line4, $line5) {
    return '<w:r><w:rPr><w:sz w:val="32"/></w:rPr>' +
           '<w:t>' + $line1 + '</w:t><w:br/>' +
           '<w:t xml:space="preserve">' + $line2 + '</w:t><w:br/>' +
           '<w:t xml:space="preserve">' + $line3 + '</w:t><w:br/>' +
           '<w:t>' + $line4 + '</w:t><w:br/>' +
           '<w:t>' + $line5 + '</w:t>' +
           '</w:r>'
}

# NB: compute each run's XML into a plain variable first, then pass that
# variable on - this PowerShell host mis-evaluates
# "Cmd arg1 arg2 (NestedCall ...)" (nested call as a positional arg), so
# we avoid that pattern entirely.

$xml = Cell-Run "82 x 85" "  8    5" "  ----" "8|    |" "2|    |"
Set-CellContent 1 1 $xml
$xml = Cell-Run "40 x 15" "  1    5" "  ----" "4|    |" "0|    |"
Set-CellContent 1 2 $xml
$xml = Cell-Run "24 x 90" "  9    0" "  ----" "2|    |" "4|    |"
Set-CellContent 1 3 $xml

$xml = Cell-Run "89 x 97" "  9    7" "  ----" "8|    |" "9|    |"
Set-CellContent 2 1 $xml
$xml = Cell-Run "95 x 15" "  1    5" "  ----" "9|    |" "5|    |"
Set-CellContent 2 2 $xml
$xml = Cell-Run "68 x 44" "  4    4" "  ----" "6|    |" "8|    |"
Set-CellContent 2 3 $xml

$xml = Cell-Run "18 x 36" "  3    6" "  ----" "1|    |" "8|    |"
Set-CellContent 3 1 $xml
$xml = Cell-Run "99 x 67" "  6    7" "  ----" "9|    |" "9|    |"
Set-CellContent 3 2 $xml
$xml = Cell-Run "33 x 46" "  4    6" "  ----" "3|    |" "3|    |"
Set-CellContent 3 3 $xml

$xml = Cell-Run "31 x 71" "  7    1" "  ----" "3|    |" "1|    |"
Set-CellContent 4 1 $xml
$xml = Cell-Run "99 x 62" "  6    2" "  ----" "9|    |" "9|    |"
Set-CellContent 4 2 $xml
$xml = Cell-Run "80 x 27" "  2    7" "  ----" "8|    |" "0|    |"
Set-CellContent 4 3 $xml

$xml = Cell-Run "59 x 87" "  8    7" "  ----" "5|    |" "9|    |"
Set-CellContent 5 1 $xml
$xml = Cell-Run "38 x 53" "  5    3" "  ----" "3|    |" "8|    |"
Set-CellContent 5 2 $xml
$xml = Cell-Run "48 x 29" "  2    9" "  ----" "4|    |" "8|    |"
Set-CellContent 5 3 $xml

Write-Host "Updated" $t.Rows.Count "x" $t.Columns.Count "lattice-multiplication table."
